$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Discipline"
$ws.Range("B2").Value = 1617.0
$ws.Range("C2").Value = "Heptathlon 100 meters hurdles"

$ws.Columns.Item(3).ColumnWidth = 32.578125
